# Auto-generated script: apply scheduled market-price refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 403739
$ws.Range("J17").Value = 403739
$ws.Range("L17").Value = 1211217
$ws.Range("N17").Value = -1211553
$ws.Range("H29").Value = 449
$ws.Range("I29").Value = 449
$ws.Range("K29").Value = 1347
$ws.Range("M29").Value = -1066
$ws.Range("H63").Value = 30246
$ws.Range("I63").Value = 30246
$ws.Range("K63").Value = 30246
$ws.Range("M63").Value = -29622
$ws.Range("H66").Value = 30246
$ws.Range("I66").Value = 30246
$ws.Range("K66").Value = 90738
$ws.Range("M66").Value = -87618
$ws.Range("H74").Value = 79436.07000000001
$ws.Range("I74").Value = 117751.664
$ws.Range("K74").Value = 117751.664
$ws.Range("M74").Value = -116815.664
$ws.Range("H77").Value = 79436.07000000001
$ws.Range("I77").Value = 117751.664
$ws.Range("K77").Value = 588758.3200000001
$ws.Range("M77").Value = -584078.3200000001
$ws.Range("H88").Value = 4100
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 4625
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 4625
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -5437
$ws.Range("H91").Value = 4100
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 4625
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 4625
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -7433
$ws.Range("H96").Value = 609.1
$ws.Range("I96").Value = 632.4
$ws.Range("J96").Value = 539.2
$ws.Range("K96").Value = 1897.2
$ws.Range("L96").Value = 1617.6
$ws.Range("M96").Value = -524.1999999999998
$ws.Range("N96").Value = -4363.6
$ws.Range("H112").Value = 78651.84
$ws.Range("J112").Value = 92857.27
$ws.Range("L112").Value = 278571.81
$ws.Range("N112").Value = -280787.81
$ws.Range("H132").Value = 20896976
$ws.Range("I132").Value = 23810624
$ws.Range("K132").Value = 71431872
$ws.Range("M132").Value = -71429342
$ws.Range("H133").Value = 79066
$ws.Range("J133").Value = 79066
$ws.Range("L133").Value = 79066
$ws.Range("N133").Value = -89186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 509995
$ws.Range("I58").Value = 999990
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 999990
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = -999560
$ws.Range("N58").Value = -20860
$ws.Range("H61").Value = 3233.9473
$ws.Range("I61").Value = 2769.9375
$ws.Range("K61").Value = 2769.9375
$ws.Range("M61").Value = -2557.9375
$ws.Range("H132").Value = 3206.7932
$ws.Range("I132").Value = 1854
$ws.Range("J132").Value = 6213
$ws.Range("K132").Value = 5562
$ws.Range("L132").Value = 18639
$ws.Range("M132").Value = -3032
$ws.Range("N132").Value = -23699
$ws.Range("H135").Value = 62850.5
$ws.Range("J135").Value = 62850.5
$ws.Range("L135").Value = 62850.5
$ws.Range("N135").Value = -72990.5
$ws.Range("H136").Value = 3233.9473
$ws.Range("I136").Value = 2769.9375
$ws.Range("K136").Value = 8309.8125
$ws.Range("M136").Value = -5759.8125
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 31294
$ws.Range("J81").Value = 31294
$ws.Range("L81").Value = 31294
$ws.Range("N81").Value = -33416
$ws.Range("H84").Value = 31294
$ws.Range("J84").Value = 31294
$ws.Range("L84").Value = 93882
$ws.Range("N84").Value = -104490
$ws.Range("H134").Value = 2274.6667
$ws.Range("I134").Value = 1963.3793
$ws.Range("K134").Value = 5890.1379
$ws.Range("M134").Value = -3355.1379

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2880.5
$ws.Range("I132").Value = 2772.1936
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 8316.5808
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -5786.5808
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 27712.723
$ws.Range("I134").Value = 16379.615
$ws.Range("J134").Value = 57178.8
$ws.Range("K134").Value = 49138.845
$ws.Range("L134").Value = 171536.4
$ws.Range("M134").Value = -46603.845
$ws.Range("N134").Value = -176606.4
$ws.Range("H137").Value = 93557.14
$ws.Range("J137").Value = 93557.14
$ws.Range("L137").Value = 93557.14
$ws.Range("N137").Value = -103757.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 93
$ws.Range("J17").Value = 182
$ws.Range("L17").Value = 546
$ws.Range("N17").Value = -884
$ws.Range("H50").Value = 1453.3077
$ws.Range("I50").Value = 321.44446
$ws.Range("J50").Value = 4000
$ws.Range("K50").Value = 964.33338
$ws.Range("L50").Value = 12000
$ws.Range("M50").Value = -483.33338
$ws.Range("N50").Value = -12962
$ws.Range("H51").Value = 1000000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1000000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3000000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3000920
$ws.Range("H53").Value = 1453.3077
$ws.Range("I53").Value = 321.44446
$ws.Range("J53").Value = 4000
$ws.Range("K53").Value = 964.33338
$ws.Range("L53").Value = 12000
$ws.Range("M53").Value = -483.33338
$ws.Range("N53").Value = -12962
$ws.Range("H55").Value = 1728.8334
$ws.Range("J55").Value = 1728.8334
$ws.Range("L55").Value = 5186.5002
$ws.Range("N55").Value = -5540.5002
$ws.Range("H94").Value = 9999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 9999
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 29997
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -31349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19702
$ws.Range("H132").Value = 3249.3684
$ws.Range("I132").Value = 3265.3125
$ws.Range("J132").Value = 3164.3333
$ws.Range("K132").Value = 9795.9375
$ws.Range("L132").Value = 9492.999899999999
$ws.Range("M132").Value = -7265.9375
$ws.Range("N132").Value = -14552.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 50000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51288
$ws.Range("H132").Value = 3143.3823
$ws.Range("I132").Value = 2775.9614
$ws.Range("K132").Value = 8327.8842
$ws.Range("M132").Value = -5797.8842
$ws.Range("H139").Value = 25325
$ws.Range("I139").Value = 25325
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 25325
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -20185
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 30000
$ws.Range("I140").Value = 30000
$ws.Range("K140").Value = 30000
$ws.Range("M140").Value = -24820
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 1560.4
$ws.Range("I100").Value = 929.5
$ws.Range("K100").Value = 1859
$ws.Range("M100").Value = -1318
$ws.Range("H126").Value = 2415.75
$ws.Range("I126").Value = 1906.6666
$ws.Range("K126").Value = 5719.9998
$ws.Range("M126").Value = -3249.9998
$ws.Range("H137").Value = 95496.5
$ws.Range("J137").Value = 95496.5
$ws.Range("L137").Value = 95496.5
$ws.Range("N137").Value = -105696.5
$ws.Range("H140").Value = 63447.25
$ws.Range("J140").Value = 63447.25
$ws.Range("L140").Value = 63447.25
$ws.Range("N140").Value = -73807.25
